$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values could be mis-parsed as numbers by Excel
# (plain single-decimal figures) must have their format locked to Text
# BEFORE the value is written, so they round-trip as the exact string,
# matching the inline-string cells in the workbook.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D11",
    "D14",
    "D16",
    "D19",
    "D20",
    "D21",
    "D22",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D39",
    "D40",
    "D42",
    "D44",
    "D46",
    "D47",
    "D48",
    "D49"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply every cell update from the diff, in sheet order.
$ws.Range("D2").Value = "44.013.13"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.358.63"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "0.676"
$ws.Range("E5").Value = "  +3.75%  "
$ws.Range("D6").Value = "238.45"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("D7").Value = "72.96"
$ws.Range("E7").Value = "  +11.73%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  +18.12%  "
$ws.Range("E10").Value = "  +6.12%  "
$ws.Range("D11").Value = "29.54"
$ws.Range("E11").Value = "  +10.52%  "
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Value = "2.709.24"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "16.85"
$ws.Range("E14").Value = "  +9.85%  "
$ws.Range("E15").Value = "  +6.96%  "
$ws.Range("D16").Value = "0.908"
$ws.Range("E16").Value = "  +7.89%  "
$ws.Range("D17").Value = "2.357.90"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "43.991.56"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "0.0000103"
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("D20").Value = "78.04"
$ws.Range("E20").Value = "  +6.00%  "
$ws.Range("D21").Value = "6.46"
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").Value = "255.36"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -3.86%  "
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  +3.62%  "
$ws.Range("D26").Value = "10.51"
$ws.Range("E26").Value = "  +6.35%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "22.46"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "1.59"
$ws.Range("E29").Value = "  +5.55%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "172.48"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "0.132"
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("E32").Value = "  +5.31%  "
$ws.Range("D33").Value = "5.20"
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("D34").Value = "0.0732"
$ws.Range("E34").Value = "  +6.55%  "
$ws.Range("D35").Value = "5.26"
$ws.Range("E35").Value = "  +5.29%  "
$ws.Range("E36").Value = "  +9.70%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "0.0270"
$ws.Range("E39").Value = "  +7.56%  "
$ws.Range("D40").Value = "19.48"
$ws.Range("E40").Value = "  +8.79%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "8.86"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("D44").Value = "0.0985"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").Value = "98.67"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "4.46"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("D48").Value = "0.182"
$ws.Range("E48").Value = "  +12.11%  "
$ws.Range("D49").Value = "2.34"
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("D50").Value = "1.439.81"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("E51").Value = "  +1.48%  "
